# Apply "final data files and prelim mapping" update to the Data Sources List sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Poverty and Median Income -> now "2018" version ---
$ws.Range("A2").Value = "Poverty and Median Income 2018"
$ws.Range("B2").Value = "2018 US Census Small Area Income and Pov Estimates (SAIPE)"
$ws.Range("C2").Value = "https://www.census.gov/data/datasets/2018/demo/saipe/2018-state-and-county.html"

# --- Row 3: Teen Pregnancy (existing row, gains a footnote in column D) ---
$ws.Range("A3").Value = "Teen Pregnancy ages 15-19"
$ws.Range("B3").Value = "SHIFT NC (Sexual Health Initiatives For Teens) "
$ws.Range("C3").Value = "http://www.shiftnc.org/data/map/northcarolina"
$ws.Range("D3").Value = "For small sample size counties= population x .031 was used for denominator in determining rate"

# --- Row 4: County Pop 2018 (new) ---
$ws.Range("A4").Value = "County Pop 2018"
$ws.Range("B4").Value = "NC Office of State Budget and Management"
$ws.Hyperlinks.Add($ws.Range("C4"), "https://demography.osbm.nc.gov/explore/dataset/county-population-estimates-standard-revised/export/?disjunctive.county&disjunctive.population&sort=-year&refine.year=July+1,+2018") | Out-Null

# --- Row 5: Child Abuse (new) ---
$ws.Range("A5").Value = "Child Abuse- substantiated abuse per 1000, 2017"
$ws.Range("B5").Value = "Annie E Casey"
$ws.Hyperlinks.Add($ws.Range("C5"), "https://datacenter.kidscount.org/") | Out-Null
$ws.Range("D5").Value = "Management Assistance for Child Welfare, Work First, and Food & Nutrition Services in North Carolina (v3.21), University of North Carolina at Chapel Hill Jordan Institute for Families. "

# --- Row 6: Children in concentrated Pov (new) ---
$ws.Range("A6").Value = "Children in concentrated Pov"
$ws.Range("B6").Value = "Casey, from ACS 2013-2017"
$ws.Range("D6").Value = "children living in tracts with 30% or more pov rate"

# --- Row 7: Elevated lead levels age1-2 (new) ---
$ws.Range("A7").Value = "Elevated lead levels age1-2"
$ws.Range("B7").Value = "Casey, from "
$ws.Range("C7").Value = " North Carolina Department of Health and Human Services, Children's Environmental Health Section: Childhood Lead Poisoning Prevention Program Surveillance Data"
$ws.Range("D7").Value = "Lead: percent of children ages 1-2 with elevated blood lead levels = 5 micrograms per deciliter (2013 and later) in North Carolina"

# --- Row 8: Juvenile Delinquency 2018 (new) ---
$ws.Range("A8").Value = "Juvenile Delinquency 2018"
$ws.Range("B8").Value = "from Casey-- North Carolina Juvenile Justice Section: Juvenile Crime Prevention Council County Data Book."
$ws.Range("D8").Value = "rate per 1,000 juveniles ages 6-15 years old of delinquent complaints received by court services offices"

# --- Row 9: No parent in workforce (new) ---
$ws.Range("A9").Value = "No parent in workforce "
$ws.Range("B9").Value = "Casey, from ACS 2013-2017"

# --- Row 10: Head of HH has no high school degree (new) ---
$ws.Range("A10").Value = "Head of HH has no high school degree"
$ws.Range("B10").Value = "Casey, from ACS 2013-2018"

# Column A got wider to fit the longer labels.
$ws.Columns.Item(1).ColumnWidth = 52.3

# Final selection left on A7 by the author before saving.
$ws.Range("A7").Select()
